{"js": "// Template city swap: \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\" -> \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\" in the dateline\n// paragraph (\"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433<spaces>@<DATE>@\"). The original run holding the\n// city name and the following run holding the padding + \"@<DATE>@\" token\n// get merged into a single run carrying the new city name (matching how\n// Word's own find/replace collapses same-formatted adjacent runs).\n\nconst body = context.document.body;\n\nconst cityWord = \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\";\nconst newCity = \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\";\nconst marker = \"@<DATE>@\";\n\n// Find the city anchor text anywhere in the body.\nconst hits = body.search(cityWord, { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find '\" + cityWord + \"' in the document body.\");\n}\n\nconst cityRange = hits.items[0];\nconst paragraph = cityRange.paragraphs.getFirst();\nparagraph.load(\"text\");\nawait context.sync();\n\n// Work out the exact run of text from the city name through the trailing\n// \"@<DATE>@\" token (inclusive of whatever whitespace padding sits between\n// them), so we don't have to hard-code the padding width.\nconst fullText = paragraph.text;\nconst startIdx = fullText.indexOf(cityWord);\nconst markerIdx = fullText.indexOf(marker, startIdx);\n\nif (startIdx === -1 || markerIdx === -1) {\n  throw new Error(\"Could not locate the city/date layout in the paragraph.\");\n}\n\nconst endIdx = markerIdx + marker.length;\nconst oldSegment = fullText.substring(startIdx, endIdx);\nconst newSegment = newCity + fullText.substring(startIdx + cityWord.length, endIdx);\n\nconst segmentHits = paragraph.search(oldSegment, { matchCase: true, matchWholeWord: false });\nsegmentHits.load(\"items\");\nawait context.sync();\n\nif (segmentHits.items.length === 0) {\n  throw new Error(\"Could not re-locate the city/date segment for replacement.\");\n}\n\n// Replacing the whole span in one shot (rather than the city word alone)\n// merges it back into a single run, just like Word's native replace does.\nsegmentHits.items[0].insertText(newSegment, \"Replace\");\nawait context.sync();\n", "ps1": "# Template city swap: \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\" -> \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\" in the dateline\n# paragraph (\"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433<spaces>@<DATE>@\"). Word's Find/Replace merges the\n# matched run back into the formatting of the surrounding text, so the run\n# that used to hold just \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\" and the following run holding the\n# padding + \"@<DATE>@\" token collapse into a single run carrying the new\n# city name - matching the target OOXML exactly.\n\n$d = $word.ActiveDocument\n\n$cityWord = \"\u0415\u043a\u0430\u0442\u0435\u0440\u0438\u043d\u0431\u0443\u0440\u0433\"\n$newCity  = \"\u041d\u043e\u0432\u043e\u0440\u043e\u0441\u0441\u0438\u0439\u0441\u043a\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = $cityWord\n$find.Replacement.Text = $newCity\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n# Positional form of Find.Execute mirrors the Word object model:\n#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#           ReplaceWith, Replace)\n# Replace:=2 is wdReplaceAll (only one match exists, so this is equivalent\n# to wdReplaceOne here too).\n$found = $find.Execute($cityWord, $true, $false, $false, $false, $false, $true, 1, $false, $newCity, 2)\n\nif (-not $found) {\n  throw \"Could not find '$cityWord' to replace with '$newCity'.\"\n}\n"}
